$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.01
$ws.Range("C4").Value = -13.318
$ws.Range("C7").Value = -13.43
$ws.Range("C8").Value = -12.672
$ws.Range("B11").Value = 6.025
$ws.Range("B12").Value = 5.729000000000001
$ws.Range("C12").Value = -13.213
$ws.Range("C14").Value = -12.048
$ws.Range("B15").Value = 6.441
$ws.Range("C22").Value = -13.318
